$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gyroscope data re-export (May 9th): two new samples inserted right after the
# header (shifting the existing 20 rows down by two) and eight new samples
# appended at the end of the series.
$data = @(
  @(2, -1.109872460365295, 1.229648113250732, 2.535529136657715),
  @(3, -0.0105194868519902, 0.2277668565511703, 1.795702934265137),
  @(4, -0.8973521590232849, 0.1004677563905716, 1.17305588722229),
  @(5, -0.6182528734207153, 1.669069647789002, -0.8818392157554626),
  @(6, 1.658749938011169, -1.585313200950623, -3.374025344848633),
  @(7, 2.416685581207275, -2.322476148605347, -5.666474342346191),
  @(8, 1.654488801956177, -1.7259281873703, -6.852646350860596),
  @(9, -3.016163110733032, -3.968842506408691, 3.157110929489136),
  @(10, -5.865745544433594, 2.865894317626953, 8.381167411804199),
  @(11, -0.1431449055671692, 12.14008617401123, 1.953362107276917),
  @(12, 4.818324089050293, -6.334262371063232, -3.150319814682007),
  @(13, 9.644504547119141, 3.98388934135437, 0.539756178855896),
  @(14, 3.38767409324646, 3.137536764144897, -4.873384952545166),
  @(15, -5.203151226043701, 0.9888983368873596, 2.413556337356567),
  @(16, -13.87813186645508, 0.5345630049705505, 0.425772875547409),
  @(17, 5.411143779754639, -7.761183738708496, -2.545782327651977),
  @(18, 2.374074935913086, -1.723265051841736, 1.63484799861908),
  @(19, 10.5457181930542, 8.709402084350586, -1.038965702056885),
  @(20, -0.5905559659004211, -1.791441917419434, 0.2574611008167267),
  @(21, -3.922703266143799, 3.348459005355835, 1.719536542892456),
  @(22, -12.29035568237305, -13.32346248626709, -5.019326210021973),
  @(23, 4.287290096282959, -8.83603572845459, -1.076782584190369),
  @(24, 1.13570511341095, -4.037551879882812, 1.869738817214966),
  @(25, 6.204765796661377, 5.688312530517578, -1.246159672737122),
  @(26, -2.400972843170166, 2.221409320831299, 0.4785034656524658),
  @(27, -4.018577098846436, 3.744737386703491, 4.383230209350586),
  @(28, -5.831124305725098, 10.37174797058106, -0.6757105588912964),
  @(29, -0.996954381465912, 15.63681697845459, -8.430303573608398),
  @(30, 5.841510772705078, -9.827264785766602, -3.761781692504883),
  @(31, 1.826529026031494, 2.439255952835083, 2.031126499176025)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
}
